$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name of the new tool being added to the shop tools list
$newTool = "Standard-Wrench_open-end"

# Rows 27-37 extend the existing list (which ran from row 2 through row 26).
# Column A keeps incrementing the running "#" counter via the same relative
# formula used by the rows above (A[n] = A[n-1] + 1).
# Column B gets the new tool name for rows 27-35; rows 36-37 are left blank
# in column B (matching the source data for this addition).
for ($r = 27; $r -le 37; $r++) {
    $ws.Cells.Item($r, 1).FormulaR1C1 = "=R[-1]C+1"
    if ($r -le 35) {
        $ws.Cells.Item($r, 2).Value = $newTool
    }
}

# Selection moves to C27, matching where the user was working next
$ws.Range("C27").Select()
